$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Add new data to row 2: category (D), title (E), and study weight (J)
$ws.Range("D2").Value = "student"
$ws.Range("E2").Value = "NAWA"
$ws.Range("J2").Value = 20

# Scroll the view so column D is the first visible column, then select K2
# as the active cell, matching the saved view state.
$win = $excel.ActiveWindow
$win.TopLeftCell = $ws.Range("D1")
$ws.Range("K2").Select() | Out-Null
